# Update "想去人数" (want-to-go count) figures for the 杭州-漫展信息 workbook.
# Same underlying event data is duplicated on the "展览" sheet and rolled up
# again on the "全部类型" sheet, so every F-column change has to be applied
# in both places. Row 12 on 展览 (row 17 on 全部类型) additionally flips from
# a numeric lowest-price to a "已售罄" (sold out) label in column G.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> new "F" (想去人数) value, for the 展览 sheet.
$exhibitionUpdates = @{
    2  = 2244
    3  = 124
    6  = 128
    7  = 51
    9  = 2636
    10 = 1629
    11 = 1651
    12 = 315
    13 = 271
    14 = 677
    15 = 845
    16 = 119
    17 = 345
    18 = 1101
    22 = 5813
    23 = 237
    24 = 1076
    25 = 121
    26 = 167
    28 = 266
    29 = 235
    30 = 49
    31 = 1068
    32 = 845
    34 = 72
    37 = 1217
    39 = 124
    42 = 135
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row 12 went from a ticket price of 218 to "已售罄" (sold out).
$sheetExhibition.Range("G12").Value = "已售罄"

# Row -> new "F" (想去人数) value, for the 全部类型 sheet (same events, offset rows).
$allTypesUpdates = @{
    2  = 2244
    5  = 124
    8  = 128
    9  = 51
    14 = 2636
    15 = 1629
    16 = 1651
    17 = 315
    18 = 271
    19 = 677
    21 = 845
    22 = 119
    23 = 345
    24 = 1101
    27 = 5813
    28 = 237
    29 = 1076
    30 = 121
    31 = 167
    33 = 266
    34 = 235
    35 = 49
    36 = 1068
    37 = 845
    39 = 72
    41 = 1217
    43 = 124
    46 = 135
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}

# Mirror the same sold-out change on row 17 of 全部类型.
$sheetAllTypes.Range("G17").Value = "已售罄"
